$d = $word.ActiveDocument

$replacements = @(
    @("86×54=", "28×62="),
    @("12×41=", "65×37="),
    @("16×53=", "34×54="),
    @("25×23=", "81×44="),
    @("73×47=", "43×22="),
    @("86×71=", "73×98="),
    @("15×60=", "22×57="),
    @("27×61=", "39×81="),
    @("12×85=", "27×12="),
    @("82×51=", "38×59="),
    @("53×15=", "14×49="),
    @("71×55=", "37×40="),
    @("69×58=", "66×82="),
    @("54×16=", "53×81="),
    @("26×85=", "74×51="),
    @("29×83=", "72×25="),
    @("75×26=", "25×72="),
    @("96×59=", "50×20="),
    @("56×76=", "74×69="),
    @("84×46=", "59×36="),
    @("55×57=", "39×14="),
    @("69×69=", "32×30="),
    @("27×91=", "91×39="),
    @("32×17=", "16×97="),
    @("40×20=", "63×86=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
